$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Add a new worksheet right after Sheet1; it becomes the active tab (Sheet2)
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Sheet2"

# Enter data in the same order the original author appears to have used
# (this reproduces the shared-string ordering seen in the target file).
$ws2.Range("B1").Value = "start"
$ws2.Range("D1").Value = "end"

$ws2.Range("A3").Value = "DSA"
$ws2.Range("A5").Value = "Web Dev"
$ws2.Range("A7").Value = "Data Science"
$ws2.Range("A9").Value = "MYSQL"

$ws2.Range("A1").Value = "Phase"
$ws2.Range("F1").Value = "days"

# Date values (stored as plain date serials, formatted as short dates)
$ws2.Range("B3").Value = 45819
$ws2.Range("D3").Value = 45930
$ws2.Range("B5").Value = 45931
$ws2.Range("D5").Value = 45961
$ws2.Range("B7").Value = 45962
$ws2.Range("D7").Value = 45976
$ws2.Range("B9").Value = 45977
$ws2.Range("D9").Value = 45991

# Apply a short-date number format to the first date cell, then copy that
# format (not value) onto the remaining date cells so they all share a
# single reused cell style instead of each creating its own style entry.
$ws2.Range("B3").NumberFormat = "mm-dd-yy"
$ws2.Range("B3").Copy()
$ws2.Range("D3").PasteSpecial(-4122)
$ws2.Range("B5").PasteSpecial(-4122)
$ws2.Range("D5").PasteSpecial(-4122)
$ws2.Range("B7").PasteSpecial(-4122)
$ws2.Range("D7").PasteSpecial(-4122)
$ws2.Range("B9").PasteSpecial(-4122)
$ws2.Range("D9").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Days values
$ws2.Range("F3").Value = 112
$ws2.Range("F5").Value = 31
$ws2.Range("F7").Value = 15
$ws2.Range("F9").Value = 15

# Column widths similar to the source workbook (values chosen so that,
# after this engine's internal character/pixel quantization, the stored
# <col> widths land as close as possible to 17, 11.5546875, 10.33203125,
# 11.5546875 respectively).
$ws2.Columns.Item(1).ColumnWidth = 16.166666666666668
$ws2.Columns.Item(2).ColumnWidth = 10.666666666666666
$ws2.Columns.Item(3).ColumnWidth = 9.5
$ws2.Columns.Item(4).ColumnWidth = 10.666666666666666

# Restore the selection/active cell on the new sheet
$ws2.Range("J16").Select()
